$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data (columns D, L, M, N, O, P, R, S) between row 4 and row 5.
# Capture current row 4 values before overwriting.
$row4 = @{
    D = $ws.Range("D4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    R = $ws.Range("R4").Value2
    S = $ws.Range("S4").Value2
}

$row5 = @{
    D = $ws.Range("D5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    R = $ws.Range("R5").Value2
    S = $ws.Range("S5").Value2
}

# Write row 5's old values into row 4
$ws.Range("D4").Value2 = $row5.D
$ws.Range("L4").Value2 = $row5.L
$ws.Range("M4").Value2 = $row5.M
$ws.Range("N4").Value2 = $row5.N
$ws.Range("O4").Value2 = $row5.O
$ws.Range("P4").Value2 = $row5.P
$ws.Range("R4").Value2 = $row5.R
$ws.Range("S4").Value2 = $row5.S

# Write row 4's old values into row 5
$ws.Range("D5").Value2 = $row4.D
$ws.Range("L5").Value2 = $row4.L
$ws.Range("M5").Value2 = $row4.M
$ws.Range("N5").Value2 = $row4.N
$ws.Range("O5").Value2 = $row4.O
$ws.Range("P5").Value2 = $row4.P
$ws.Range("R5").Value2 = $row4.R
$ws.Range("S5").Value2 = $row4.S
